$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of bird observation data appended to the sheet (rows 62-64)
$rows = @{}

$rows["62"] = [ordered]@{
    "A" = @{ Type = "n"; Value = 131273902 }
    "B" = @{ Type = "n"; Value = 57884 }
    "D" = @{ Type = "s"; Value = "NT" }
    "E" = @{ Type = "n"; Value = 100109 }
    "F" = @{ Type = "s"; Value = "Tretåig hackspett" }
    "G" = @{ Type = "s"; Value = "Picoides tridactylus" }
    "H" = @{ Type = "s"; Value = "(Linnaeus, 1758)" }
    "I" = @{ Type = "s"; Value = "" }
    "K" = @{ Type = "s"; Value = "" }
    "L" = @{ Type = "s"; Value = "" }
    "M" = @{ Type = "s"; Value = "äldre spår" }
    "N" = @{ Type = "s"; Value = "" }
    "P" = @{ Type = "s"; Value = "Sims fäbodar, Dlr" }
    "Q" = @{ Type = "n"; Value = 515081 }
    "R" = @{ Type = "n"; Value = 6704854 }
    "S" = @{ Type = "n"; Value = 50 }
    "T" = @{ Type = "s"; Value = "Dalarna" }
    "U" = @{ Type = "s"; Value = "Borlänge" }
    "V" = @{ Type = "s"; Value = "Dalarna" }
    "W" = @{ Type = "s"; Value = "Stora Tuna" }
    "Y" = @{ Type = "d"; Value = "2026-02-23" }
    "AA" = @{ Type = "d"; Value = "2026-02-23" }
    "AC" = @{ Type = "s"; Value = "Ringhack på tall." }
    "AD" = @{ Type = "b"; Value = $false }
    "AE" = @{ Type = "b"; Value = $false }
    "AG" = @{ Type = "b"; Value = $false }
    "AT" = @{ Type = "s"; Value = "" }
    "AW" = @{ Type = "s"; Value = "Anna-Lena Thommson" }
    "AX" = @{ Type = "s"; Value = "Anna-Lena Thommson" }
    "AY" = @{ Type = "s"; Value = "" }
}

$rows["63"] = [ordered]@{
    "A" = @{ Type = "n"; Value = 131273875 }
    "B" = @{ Type = "n"; Value = 57884 }
    "D" = @{ Type = "s"; Value = "NT" }
    "E" = @{ Type = "n"; Value = 100109 }
    "F" = @{ Type = "s"; Value = "Tretåig hackspett" }
    "G" = @{ Type = "s"; Value = "Picoides tridactylus" }
    "H" = @{ Type = "s"; Value = "(Linnaeus, 1758)" }
    "I" = @{ Type = "s"; Value = "" }
    "K" = @{ Type = "s"; Value = "" }
    "L" = @{ Type = "s"; Value = "" }
    "M" = @{ Type = "s"; Value = "färska spår" }
    "N" = @{ Type = "s"; Value = "" }
    "P" = @{ Type = "s"; Value = "Sims fäbodar, Dlr" }
    "Q" = @{ Type = "n"; Value = 515025 }
    "R" = @{ Type = "n"; Value = 6704768 }
    "S" = @{ Type = "n"; Value = 50 }
    "T" = @{ Type = "s"; Value = "Dalarna" }
    "U" = @{ Type = "s"; Value = "Borlänge" }
    "V" = @{ Type = "s"; Value = "Dalarna" }
    "W" = @{ Type = "s"; Value = "Stora Tuna" }
    "Y" = @{ Type = "d"; Value = "2026-02-23" }
    "AA" = @{ Type = "d"; Value = "2026-02-23" }
    "AC" = @{ Type = "s"; Value = "Ringhack på tall." }
    "AD" = @{ Type = "b"; Value = $false }
    "AE" = @{ Type = "b"; Value = $false }
    "AG" = @{ Type = "b"; Value = $false }
    "AT" = @{ Type = "s"; Value = "" }
    "AW" = @{ Type = "s"; Value = "Anna-Lena Thommson" }
    "AX" = @{ Type = "s"; Value = "Anna-Lena Thommson" }
    "AY" = @{ Type = "s"; Value = "" }
}

$rows["64"] = [ordered]@{
    "A" = @{ Type = "n"; Value = 131273771 }
    "B" = @{ Type = "n"; Value = 57073 }
    "D" = @{ Type = "s"; Value = "LC" }
    "E" = @{ Type = "n"; Value = 100138 }
    "F" = @{ Type = "s"; Value = "Tjäder" }
    "G" = @{ Type = "s"; Value = "Tetrao urogallus" }
    "H" = @{ Type = "s"; Value = "Linnaeus, 1758" }
    "I" = @{ Type = "s"; Value = "" }
    "K" = @{ Type = "s"; Value = "" }
    "L" = @{ Type = "s"; Value = "" }
    "M" = @{ Type = "s"; Value = "färsk spillning" }
    "N" = @{ Type = "s"; Value = "" }
    "P" = @{ Type = "s"; Value = "Sims fäbodar, Dlr" }
    "Q" = @{ Type = "n"; Value = 515338 }
    "R" = @{ Type = "n"; Value = 6704939 }
    "S" = @{ Type = "n"; Value = 50 }
    "T" = @{ Type = "s"; Value = "Dalarna" }
    "U" = @{ Type = "s"; Value = "Borlänge" }
    "V" = @{ Type = "s"; Value = "Dalarna" }
    "W" = @{ Type = "s"; Value = "Stora Tuna" }
    "Y" = @{ Type = "d"; Value = "2026-02-23" }
    "AA" = @{ Type = "d"; Value = "2026-02-23" }
    "AD" = @{ Type = "b"; Value = $false }
    "AE" = @{ Type = "b"; Value = $false }
    "AG" = @{ Type = "b"; Value = $false }
    "AT" = @{ Type = "s"; Value = "" }
    "AW" = @{ Type = "s"; Value = "Anna-Lena Thommson" }
    "AX" = @{ Type = "s"; Value = "Anna-Lena Thommson" }
    "AY" = @{ Type = "s"; Value = "" }
}

# Cells holding date-like text must be forced to Text format first,
# otherwise Excel auto-converts the "YYYY-MM-DD" string into a date serial.
$dateCells = New-Object System.Collections.ArrayList

# Rows are written in ascending order so the sheet grows top-to-bottom naturally.
$rowNumbers = $rows.Keys | Sort-Object { [int]$_ }
foreach ($rownum in $rowNumbers) {
    $rowdata = $rows[$rownum]
    foreach ($col in $rowdata.Keys) {
        $cellinfo = $rowdata[$col]
        $addr = "$col$rownum"
        $range = $ws.Range($addr)
        if ($cellinfo.Type -eq "d") {
            $range.NumberFormat = "@"
            $range.Value = $cellinfo.Value
            [void]$dateCells.Add($addr)
        } else {
            $range.Value = $cellinfo.Value
        }
    }
}

# Strip the temporary Text number format back off so the cells carry no
# explicit style, while the stored value remains the literal date string.
foreach ($addr in $dateCells) {
    $ws.Range($addr).ClearFormats()
}

